$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 96, shifting existing rows 96-184 down to 97-185.
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new record.
$ws.Range("A96").Value = 5
$ws.Range("B96").Value = "Macroferia Regional de Talca"
$ws.Range("C96").Value = "Maule"
$ws.Range("D96").Value = 44874
$ws.Range("E96").Value = 7
$ws.Range("F96").Value = 100112031
$ws.Range("G96").Value = "Poroto verde"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = 2000
$ws.Range("N96").Value = "$/kilo"
$ws.Range("O96").Value = "Región del Maule"
$ws.Range("P96").Value = 2000
$ws.Range("Q96").Value = 1
$ws.Range("R96").Value = "Hortaliza"
